$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("This message appears when a person logs in.")
$para = $rng.Paragraphs(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="07A5C4E9" w14:textId="77777777" w:rsidR="00E40122" w:rsidRDefault="00E40122" w:rsidP="00DF738E"><w:r><w:t>This message appears when a person logs in.</w:t></w:r></w:p>'
$para.Range.InsertXML($xml)
